$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.671.36"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "3.516.88"
$ws.Range("E3").Value = "  -1.85%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "586.61"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.74%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "132.91"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D7").Value = "3.517.63"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -0.46%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "4.115.09"
$ws.Range("E13").Value = "  -1.87%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.85"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000180"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "3.516.75"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "64.624.89"
$ws.Range("E18").Value = "  -0.20%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.96"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.88%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.26"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.71"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.60%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "391.37"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "3.657.95"
$ws.Range("E24").Value = "  -1.78%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.11"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E27").Value = "  -4.08%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.75%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.46"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -7.99%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  -4.25%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.23"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("D33").Value = "3.520.79"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.66%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.146"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  -0.10%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "171.26"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0811"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.16%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.815"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.03%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "26.43"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  -0.04%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.05"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("E46").Value = "  -2.85%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").Value = "2.474.47"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0269"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
